$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 497, shifting existing data (rows 497-528) down to 498-529.
$ws.Rows.Item(497).Insert()

# Copy the date number format (style) used by column D in the surrounding rows onto the new D497 cell.
$ws.Range("D496").Copy()
$ws.Range("D497").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the newly inserted row 497 with its values.
$ws.Cells.Item(497, 1).Value = 10
$ws.Cells.Item(497, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(497, 3).Value = "La Araucanía"
$ws.Cells.Item(497, 4).Value = 44931
$ws.Cells.Item(497, 5).Value = 9
$ws.Cells.Item(497, 6).Value = 100112008
$ws.Cells.Item(497, 7).Value = "Coliflor"
$ws.Cells.Item(497, 8).Value = "Sin especificar"
$ws.Cells.Item(497, 9).Value = "Primera"
$ws.Cells.Item(497, 10).Value = 500
$ws.Cells.Item(497, 11).Value = 1200
$ws.Cells.Item(497, 12).Value = 1200
$ws.Cells.Item(497, 13).Value = 1200
$ws.Cells.Item(497, 14).Value = "`$/unidad"
$ws.Cells.Item(497, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(497, 16).Value = 1200
$ws.Cells.Item(497, 17).Value = 1
$ws.Cells.Item(497, 18).Value = "Hortaliza"
